# EPBDS-6830 Implement circular datatype dependencies support like in Java.
# This script mutates the test workbook:
#  - On the "Test" sheet: replaces the hard-coded "enum bit" numbers
#    (1000/2000/4000/.../64000) with their simplified powers of two
#    (1/2/4/8/16/32/64), and renames the leftover "MyTXT" label to "message".
#  - On the "Data" sheet: the same numeric simplification plus the same
#    "MyTXT" -> "message" rename.
# Once no cell references the string "MyTXT" any more, it naturally drops
# out of the workbook's shared string table when the file is saved.

$wb = $excel.ActiveWorkbook

$testSheet = $wb.Worksheets.Item("Test")
$dataSheet = $wb.Worksheets.Item("Data")

# ---- Test sheet -----------------------------------------------------
$testSheet.Range("J7").Value = 1
$testSheet.Range("K7").Value = 2
$testSheet.Range("L7").Value = 4
$testSheet.Range("M7").Value = 8
$testSheet.Range("N7").Value = 16
$testSheet.Range("O7").Value = 32

$testSheet.Range("K8").Value = 2
$testSheet.Range("L8").Value = 4
$testSheet.Range("M8").Value = 8
$testSheet.Range("N8").Value = 16
$testSheet.Range("O8").Value = 32
$testSheet.Range("P8").Value = 64

$testSheet.Range("K16").Value = 1
$testSheet.Range("L16").Value = 2
$testSheet.Range("M16").Value = 4
$testSheet.Range("N16").Value = 8
$testSheet.Range("O16").Value = 16
$testSheet.Range("P16").Value = 32
$testSheet.Range("R16").Value = "message"

$testSheet.Range("L17").Value = 2
$testSheet.Range("M17").Value = 4
$testSheet.Range("N17").Value = 8
$testSheet.Range("O17").Value = 16
$testSheet.Range("P17").Value = 32
$testSheet.Range("Q17").Value = 64
$testSheet.Range("R17").Value = "message"

# ---- Data sheet -------------------------------------------------------
$dataSheet.Range("C11").Value = 1
$dataSheet.Range("D11").Value = 2
$dataSheet.Range("E11").Value = 4
$dataSheet.Range("F11").Value = 8
$dataSheet.Range("G11").Value = 16
$dataSheet.Range("H11").Value = 32

$dataSheet.Range("O11").Value = 1
$dataSheet.Range("P11").Value = 2
$dataSheet.Range("Q11").Value = 4
$dataSheet.Range("R11").Value = 8
$dataSheet.Range("S11").Value = 16
$dataSheet.Range("T11").Value = 32
$dataSheet.Range("V11").Value = "message"

$dataSheet.Range("D12").Value = 2
$dataSheet.Range("E12").Value = 4
$dataSheet.Range("F12").Value = 8
$dataSheet.Range("G12").Value = 16
$dataSheet.Range("H12").Value = 32
$dataSheet.Range("I12").Value = 64

$dataSheet.Range("P12").Value = 2
$dataSheet.Range("Q12").Value = 4
$dataSheet.Range("R12").Value = 8
$dataSheet.Range("S12").Value = 16
$dataSheet.Range("T12").Value = 32
$dataSheet.Range("U12").Value = 64
$dataSheet.Range("V12").Value = "message"

$dataSheet.Range("C37").Value = 1
$dataSheet.Range("D37").Value = 2
$dataSheet.Range("E37").Value = 4
$dataSheet.Range("F37").Value = 8
$dataSheet.Range("G37").Value = 16
$dataSheet.Range("H37").Value = 32

$dataSheet.Range("D38").Value = 2
$dataSheet.Range("E38").Value = 4
$dataSheet.Range("F38").Value = 8
$dataSheet.Range("G38").Value = 16
$dataSheet.Range("H38").Value = 32
$dataSheet.Range("I38").Value = 64

$wb.Save()
